$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("CZ2").Value = 0.004121571094013995
$ws.Range("DA2").Value = 0.0003144382496725302
$ws.Range("DB2").Value = 0.0003416011628651941
$ws.Range("DC2").Value = 0.000001052918379933622
$ws.Range("DD2").Value = 0.008066191254859696
$ws.Range("DE2").Value = 0.0006562996704702154
$ws.Range("DF2").Value = 0.000003095910523986615
$ws.Range("DG2").Value = 0.0006638803260742153
$ws.Range("DH2").Value = 0.01203026486521681
$ws.Range("DN2").Value = 0.003196148776567774
$ws.Range("DO2").Value = 0.0001395886312723529
$ws.Range("DP2").Value = 0.004121571094013995
$ws.Range("DQ2").Value = 0.0003144382496725302
$ws.Range("DR2").Value = 0.003658292532094398
$ws.Range("DS2").Value = 0.0002739988145328696
$ws.Range("DT2").Value = 0.008111074532716082
$ws.Range("DU2").Value = 0.0006638803260742153
$ws.Range("DV2").Value = 0.01582596971053386
$ws.Range("DW2").Value = 0.0007588546929428664
$ws.Range("DX2").Value = 10.26194040763103
$ws.Range("DY2").Value = 0.01582596971053386
$ws.Range("DZ2").Value = 37.27671427850061
$ws.Range("EA2").Value = 0.02033539122703944
$ws.Range("EB2").Value = 0.3736013441107552
$ws.Range("EC2").Value = 0.02644687621329951
$ws.Range("ED2").Value = 0.0406707824540792
$ws.Range("EE2").Value = 0.007938000104536963
$ws.Range("EF2").Value = 0.0001402195111432296
$ws.Range("EG2").Value = 24.0898168839625
$ws.Range("EH2").Value = 0.007938000104536963
$ws.Range("EI2").Value = 60.12322451458549
$ws.Range("EJ2").Value = 0.002761838881079758
$ws.Range("EL2").Value = 0.1002807771246959
$ws.Range("EM2").Value = 0.02209513620528016
$ws.Range("EN2").Value = 0.005523677762159264
$ws.Range("FD2").Value = 0.1089431235373258
$ws.Range("FE2").Value = 0.08354097332795805
$ws.Range("FF2").Value = 0.2206193653076095
$ws.Range("FG2").Value = 0.0001281248056977315
$ws.Range("FH2").Value = 0.02267008832946689
$ws.Range("FI2").Value = 0.03604932200671356
$ws.Range("FJ2").Value = 0.03656020773586825
$ws.Range("FK2").Value = 0.0001264591013608934
$ws.Range("FN2").Value = 0.9847319278346619
$ws.Range("FP2").Value = 318.8004030262446
$ws.Range("FV2").Value = 318.8004030262446
$ws.Range("FW2").Value = 0.0000002602965943640123
$ws.Range("FX2").Value = 0.0000002887859710303373
$ws.Range("FY2").Value = 0.0000002318077307103913
$ws.Range("CZ3").Value = 0.003898126331607524
$ws.Range("DA3").Value = 0.0005436650694396729
$ws.Range("DB3").Value = 0.0007002384953520575
$ws.Range("DC3").Value = 0.0000009738064076524554
$ws.Range("DD3").Value = 0.007665036205855074
$ws.Range("DE3").Value = 0.001379276957156411
$ws.Range("DF3").Value = 0.000003966936651187698
$ws.Range("DG3").Value = 0.001442328582060869
$ws.Range("DH3").Value = 0.009025217197084384
$ws.Range("DN3").Value = 0.00303719504256277
$ws.Range("DO3").Value = 0.0001389676585686582
$ws.Range("DP3").Value = 0.003898126331607524
$ws.Range("DQ3").Value = 0.0005436650694396729
$ws.Range("DR3").Value = 0.003506524124719863
$ws.Range("DS3").Value = 0.0005400226476868536
$ws.Range("DT3").Value = 0.007684615965818736
$ws.Range("DU3").Value = 0.001442328582060869
$ws.Range("DV3").Value = 0.009662549737418471
$ws.Range("DW3").Value = 0.001591873738949394
$ws.Range("DX3").Value = 7.932053356119781
$ws.Range("DY3").Value = 0.009653159104427272
$ws.Range("DZ3").Value = 30.64367133114799
$ws.Range("EA3").Value = 0.01846395916693758
$ws.Range("EB3").Value = 1.184833557582112
$ws.Range("EC3").Value = 0.01491220130670441
$ws.Range("ED3").Value = 0.03692791833387508
$ws.Range("EE3").Value = 0.002024293746508091
$ws.Range("EF3").Value = 0.00001880011025546369
$ws.Range("EG3").Value = 8.054951021302871
$ws.Range("EH3").Value = 0.002024293746493335
$ws.Range("EI3").Value = 20.40526935626539
$ws.Range("EJ3").Value = 0.002285597454323184
$ws.Range("EL3").Value = 0.1005832576141919
$ws.Range("EM3").Value = 0.0042787150575013
$ws.Range("EN3").Value = 0.00457119490864614
$ws.Range("ER3").Value = 0.0515846589239724
$ws.Range("ES3").Value = 5.542381561145443
$ws.Range("ET3").Value = 0.02315520927181866
$ws.Range("FD3").Value = 0.01574314760709608
$ws.Range("FE3").Value = 0.01303335717835754
$ws.Range("FF3").Value = 0.02993970484831542
$ws.Range("FG3").Value = 0.00006613373658423826
$ws.Range("FH3").Value = 0.005602193457839931
$ws.Range("FI3").Value = 0.008057642605724875
$ws.Range("FJ3").Value = 0.009267880971006033
$ws.Range("FK3").Value = 0.000004767334751171718
$ws.Range("FP3").Value = 231.7558799227948
$ws.Range("FV3").Value = 231.7558799227948
$ws.Range("FW3").Value = 0.0000001892258155404329
$ws.Range("FX3").Value = 0.0000001511430976542597
$ws.Range("FY3").Value = 0.0000002273086782885878
$ws.Range("CZ4").Value = 0.008633122383273701
$ws.Range("DA4").Value = 0.00005377537467975624
$ws.Range("DB4").Value = 0.00008522464012497288
$ws.Range("DC4").Value = 0.000002165665294197172
$ws.Range("DD4").Value = 0.01696085881841999
$ws.Range("DE4").Value = 0.0001704666858567483
$ws.Range("DF4").Value = 0.000002774967990583252
$ws.Range("DG4").Value = 0.0001851512566225278
$ws.Range("DH4").Value = 0.01801682715147323
$ws.Range("DN4").Value = 0.006720573124179088
$ws.Range("DO4").Value = 0.00008026288111173488
$ws.Range("DP4").Value = 0.008633122383273701
$ws.Range("DQ4").Value = 0.00005377537467975624
$ws.Range("DR4").Value = 0.007746858164322312
$ws.Range("DS4").Value = 0.00006107509112899551
$ws.Range("DT4").Value = 0.0169828414145919
$ws.Range("DU4").Value = 0.0001851512566225278
$ws.Range("DV4").Value = 0.02623274408308158
$ws.Range("DW4").Value = 0.0001586544485687383
$ws.Range("DX4").Value = 38.64868055868264
$ws.Range("DY4").Value = 0.02623274408299633
$ws.Range("DZ4").Value = 93.84580696192594
$ws.Range("EA4").Value = 0.01003618555535543
$ws.Range("EB4").Value = 0.01610203031423975
$ws.Range("EC4").Value = 0.008112346375556006
$ws.Range("ED4").Value = 0.02007237111071136
$ws.Range("EE4").Value = 0.009722221465514234
$ws.Range("EF4").Value = 0.0003023125584319004
$ws.Range("EG4").Value = 80.81872789310843
$ws.Range("EH4").Value = 0.009722221465514234
$ws.Range("EI4").Value = 173.6694211961059
$ws.Range("EJ4").Value = 0.005697831658283635
$ws.Range("EL4").Value = 0.3168704501958851
$ws.Range("EM4").Value = 0.006202739751020928
$ws.Range("EN4").Value = 0.01139566331656727
$ws.Range("EO4").Value = 0.06966309548848572
$ws.Range("EP4").Value = 21.00045867682866
$ws.Range("EQ4").Value = 0.03046794478687258
$ws.Range("ER4").Value = 0.03264958566623857
$ws.Range("ES4").Value = 21.11338136863559
$ws.Range("ET4").Value = 0.08988017537955771
$ws.Range("EU4").Value = 0.01617142175980815
$ws.Range("EV4").Value = 4.210346075723318
$ws.Range("EW4").Value = 0.003158900176068705
$ws.Range("FD4").Value = 0.02330889875037172
$ws.Range("FE4").Value = 0.01489791352539781
$ws.Range("FF4").Value = 0.04714792974045415
$ws.Range("FG4").Value = 0.0001327761391740501
$ws.Range("FH4").Value = 0.01536774808297608
$ws.Range("FI4").Value = 0.02260295803459265
$ws.Range("FJ4").Value = 0.02442674130379713
$ws.Range("FK4").Value = 0.00005059928622059757
$ws.Range("FP4").Value = 269.1238309099282
$ws.Range("FV4").Value = 269.1238309099282
$ws.Range("FW4").Value = 0.0000002197362863264472
$ws.Range("FX4").Value = 0.0000002478240697835068
$ws.Range("FY4").Value = 0.0000001916488049760567
$ws.Range("CZ5").Value = 0.008623519319753378
$ws.Range("DA5").Value = 0.001008689144248397
$ws.Range("DB5").Value = 0.001127815431212898
$ws.Range("DC5").Value = 0.000002177860760942002
$ws.Range("DD5").Value = 0.01691807483611963
$ws.Range("DE5").Value = 0.002254215106701593
$ws.Range("DF5").Value = 0.000004135608388685828
$ws.Range("DG5").Value = 0.00226601662292915
$ws.Range("DH5").Value = 0.007816067367347333
$ws.Range("DN5").Value = 0.00670362039291657
$ws.Range("DO5").Value = 0.0002920617835899507
$ws.Range("DP5").Value = 0.008623519319753378
$ws.Range("DQ5").Value = 0.001008689144248397
$ws.Range("DR5").Value = 0.007707357191337734
$ws.Range("DS5").Value = 0.000994240163399574
$ws.Range("DT5").Value = 0.01698254536139131
$ws.Range("DU5").Value = 0.00226601662292915
$ws.Range("DV5").Value = 0.03721516730204787
$ws.Range("DW5").Value = 0.001708382298420771
$ws.Range("DX5").Value = 132.9050131615059
$ws.Range("DY5").Value = 0.03721516730204787
$ws.Range("DZ5").Value = 369.2303515964784
$ws.Range("EA5").Value = 0.007001667935537704
$ws.Range("EB5").Value = 0.5967494502662745
$ws.Range("EC5").Value = 0.03334176427172011
$ws.Range("ED5").Value = 0.0140033358710753
$ws.Range("EE5").Value = 0.0207004448332303
$ws.Range("EF5").Value = 0.001692419490779649
$ws.Range("EG5").Value = 256.4716665433589
$ws.Range("EH5").Value = 0.02072324478573048
$ws.Range("EI5").Value = 594.4847793625852
$ws.Range("EJ5").Value = 0.003432779748085921
$ws.Range("EL5").Value = 0.1832455330296621
$ws.Range("EM5").Value = 0.01397411128870678
$ws.Range("EN5").Value = 0.006865559496171732
$ws.Range("ER5").Value = 0.01577966434030698
$ws.Range("ES5").Value = 11.98465134164255
$ws.Range("ET5").Value = 0.03991293541689809
$ws.Range("FD5").Value = 0.01959499579119182
$ws.Range("FE5").Value = 0.02115636098042637
$ws.Range("FF5").Value = 0.03818742823538833
$ws.Range("FG5").Value = 0.001804249583834256
$ws.Range("FH5").Value = 0.009273102205938688
$ws.Range("FI5").Value = 0.015324679381836
$ws.Range("FJ5").Value = 0.0154517210410927
$ws.Range("FK5").Value = 0.0007298877547455057
$ws.Range("FP5").Value = 218.7939393099986
$ws.Range("FV5").Value = 218.7939393099986
$ws.Range("FW5").Value = 0.0000002895706280777925
$ws.Range("FX5").Value = 0.0000003281420848418987
$ws.Range("FY5").Value = 0.0000002509994217293434
$ws.Range("CZ6").Value = 0.00313839773005612
$ws.Range("DA6").Value = 0.000086997684489125
$ws.Range("DB6").Value = 0.000103536001165768
$ws.Range("DC6").Value = 0.0000007974811352240221
$ws.Range("DD6").Value = 0.006149107031357073
$ws.Range("DE6").Value = 0.0001965583359755892
$ws.Range("DF6").Value = 0.000002473346072500333
$ws.Range("DG6").Value = 0.000198029563477719
$ws.Range("DH6").Value = 0.03942618726714858
$ws.Range("DN6").Value = 0.002436523049628237
$ws.Range("DO6").Value = 0.00004707667461508285
$ws.Range("DP6").Value = 0.00313839773005612
$ws.Range("DQ6").Value = 0.000086997684489125
$ws.Range("DR6").Value = 0.00279468786112772
$ws.Range("DS6").Value = 0.00007580422235349371
$ws.Range("DT6").Value = 0.006184902148396849
$ws.Range("DU6").Value = 0.0001980295634777183
$ws.Range("DV6").Value = 0.01546795107595128
$ws.Range("DW6").Value = 0.0002317488688330182
$ws.Range("DX6").Value = 35.81175557957203
$ws.Range("DY6").Value = 0.01546795107595128
$ws.Range("DZ6").Value = 110.0510971335006
$ws.Range("EA6").Value = 0.009151238075322976
$ws.Range("EB6").Value = 0.3762553102987861
$ws.Range("EC6").Value = 0.01999253977346099
$ws.Range("ED6").Value = 0.01830247615064589
$ws.Range("EE6").Value = 0.009362035925299618
$ws.Range("EF6").Value = 0.0001405402809915878
$ws.Range("EG6").Value = 84.74218874542244
$ws.Range("EH6").Value = 0.009362035925428435
$ws.Range("EI6").Value = 203.5671299610883
$ws.Range("EJ6").Value = 0.005150035946231741
$ws.Range("EL6").Value = 0.06253602222536761
$ws.Range("EM6").Value = 0.00812898246043466
$ws.Range("EN6").Value = 0.01030007189246353
$ws.Range("ER6").Value = 0.004437404922780453
$ws.Range("ES6").Value = 9.361221150633545
$ws.Range("ET6").Value = 0.132462260989633
$ws.Range("EX6").Value = 1.691756401039812
$ws.Range("EY6").Value = 35.9793962520101
$ws.Range("EZ6").Value = 1.605020890486599
$ws.Range("FD6").Value = 0.02848115603920553
$ws.Range("FE6").Value = 0.01189271839528398
$ws.Range("FF6").Value = 0.05756832060630127
$ws.Range("FG6").Value = 0.0001829473189592143
$ws.Range("FH6").Value = 0.001809710100015624
$ws.Range("FI6").Value = 0.01925226189805627
$ws.Range("FJ6").Value = 0.003420848607510259
$ws.Range("FK6").Value = 0.00009604395130429959
$ws.Range("FP6").Value = 259.3493023049413
$ws.Range("FV6").Value = 259.3493023049413
$ws.Range("FW6").Value = 0.0000003432450670697337
$ws.Range("FX6").Value = 0.0000003193760193205641
$ws.Range("FY6").Value = 0.0000003671147461962039
$ws.Range("CZ10").Value = 0.002944435175931269
$ws.Range("DA10").Value = 0.000212858836560792
$ws.Range("DB10").Value = 0.000273105417757433
$ws.Range("DC10").Value = 0.000002650830837543942
$ws.Range("DD10").Value = 0.007430936081823092
$ws.Range("DE10").Value = 0.0006892421743689242
$ws.Range("DF10").Value = 0.0000004515987150456344
$ws.Range("DG10").Value = 0.0006957109852536883
$ws.Range("DH10").Value = 0.001069044967649705
$ws.Range("DN10").Value = 0.002944435175931269
$ws.Range("DO10").Value = 0.000212858836560792
$ws.Range("DP10").Value = 0.002939983176574723
$ws.Range("DQ10").Value = 0.0002462028080304334
$ws.Range("DR10").Value = 0.005611446657419663
$ws.Range("DS10").Value = 0.0004896797921157245
$ws.Range("DT10").Value = 0.007459142723293881
$ws.Range("DU10").Value = 0.0006957109852536883
$ws.Range("DV10").Value = 0.008054865258021768
$ws.Range("DW10").Value = 0.0007076727748828312
$ws.Range("DX10").Value = 9.522545037403352
$ws.Range("DY10").Value = 0.008054865258021768
$ws.Range("DZ10").Value = 14.07012067614134
$ws.Range("EA10").Value = 0.0194013863993258
$ws.Range("EB10").Value = 0.121181512289704
$ws.Range("EC10").Value = 0.1119602641518845
$ws.Range("ED10").Value = 0.0388027727986516
$ws.Range("EE10").Value = 0.0005957225347026574
$ws.Range("EF10").Value = 0.0001737156748697879
$ws.Range("EG10").Value = 11.00501273265161
$ws.Range("EH10").Value = 0.0005957225347026574
$ws.Range("EI10").Value = 11.10811209393097
$ws.Range("EJ10").Value = 0.01183208200932955
$ws.Range("EL10").Value = 0.06840653180632149
$ws.Range("EM10").Value = 0.05465860951598416
$ws.Range("EN10").Value = 0.02366416401865911
$ws.Range("FD10").Value = 0.03714288941969795
$ws.Range("FE10").Value = 0.04858505470497238
$ws.Range("FF10").Value = 0.03308955508327872
$ws.Range("FG10").Value = 0.0007693052004259556
$ws.Range("FH10").Value = 0.001405384427117133
$ws.Range("FI10").Value = 0.02053751632939136
$ws.Range("FJ10").Value = 0.004728223840575032
$ws.Range("FK10").Value = 0.0002628959101822814
$ws.Range("FP10").Value = 140.5794132459365
$ws.Range("FV10").Value = 140.5794132459365
$ws.Range("FW10").Value = 0.0000001147813558907516
$ws.Range("FX10").Value = 0.0000001181359191710326
$ws.Range("FY10").Value = 0.0000001114267925569197
$ws.Range("CZ34").Value = 0.001905450080163094
$ws.Range("DA34").Value = 0.0005577006691517093
$ws.Range("DB34").Value = 0.0007336438446951432
$ws.Range("DC34").Value = 0.00000303413460406939
$ws.Range("DD34").Value = 0.004808826449459376
$ws.Range("DE34").Value = 0.001851513173503015
$ws.Range("DF34").Value = 0.0000004303726260355536
$ws.Range("DG34").Value = 0.001858474809303476
$ws.Range("DH34").Value = 0.06434000593438227
$ws.Range("DN34").Value = 0.001905450080163094
$ws.Range("DO34").Value = 0.0005577006691517093
$ws.Range("DP34").Value = 0.001923997973759467
$ws.Range("DQ34").Value = 0.0006545936307662285
$ws.Range("DR34").Value = 0.003505268867990076
$ws.Range("DS34").Value = 0.00131818159885301
$ws.Range("DT34").Value = 0.004838770837617376
$ws.Range("DU34").Value = 0.001858474809303476
$ws.Range("DV34").Value = 0.004817561442466234
$ws.Range("DW34").Value = 0.001604104821853441
$ws.Range("DX34").Value = 32.56798300279781
$ws.Range("DY34").Value = 0.004817561442260395
$ws.Range("DZ34").Value = 41.20542208866173
$ws.Range("EA34").Value = 0.009531984300892376
$ws.Range("EB34").Value = 0.6469276446995978
$ws.Range("EC34").Value = 0.03392545124800647
$ws.Range("ED34").Value = 0.01906396860178466
$ws.Range("EE34").Value = 0.002389435620152463
$ws.Range("EF34").Value = 0.003934177073118203
$ws.Range("EG34").Value = 65.9955711304164
$ws.Range("EH34").Value = 0.002359150740863949
$ws.Range("EI34").Value = 72.48635749835513
$ws.Range("EJ34").Value = 0.006612690680070731
$ws.Range("EL34").Value = 0.80274450619363
$ws.Range("EM34").Value = 0.01279819516124776
$ws.Range("EN34").Value = 0.01322538136014166
$ws.Range("EO34").Value = 0.0914035708103241
$ws.Range("EP34").Value = 7.872968064902124
$ws.Range("EQ34").Value = 0.2402978115426339
$ws.Range("ER34").Value = 0.05084695122519476
$ws.Range("ES34").Value = 6.17381440443843
$ws.Range("ET34").Value = 0.126488734252547
$ws.Range("FD34").Value = 0.01557846321835421
$ws.Range("FE34").Value = 0.02016372164802749
$ws.Range("FF34").Value = 0.02045659527001582
$ws.Range("FG34").Value = 0.0002771051074289666
$ws.Range("FH34").Value = 0.004978285176939423
$ws.Range("FI34").Value = 0.01041880326115391
$ws.Range("FJ34").Value = 0.004623255966375513
$ws.Range("FK34").Value = 0.0001544613887772369
$ws.Range("FP34").Value = 265.772288041678
$ws.Range("FQ34").Value = 0.03185698399946534
$ws.Range("FR34").Value = 6.486215243485586
$ws.Range("FS34").Value = 6.005586200685641
$ws.Range("FT34").Value = 0.004861946222160328
$ws.Range("FU34").Value = 0.002008632512260287
$ws.Range("FV34").Value = 265.772288041678
$ws.Range("FW34").Value = 0.0000003517457962489164
$ws.Range("FX34").Value = 0.0000003474410063831506
$ws.Range("FY34").Value = 0.0000003560514935761004
$ws.Range("CZ35").Value = 0.006136435233163064
$ws.Range("DA35").Value = 0.01015222382802648
$ws.Range("DB35").Value = 0.01075710415789199
$ws.Range("DC35").Value = 0.00001404147325705511
$ws.Range("DD35").Value = 0.0154866571220783
$ws.Range("DE35").Value = 0.02714794133563512
$ws.Range("DF35").Value = 0.000002589797094548286
$ws.Range("DG35").Value = 0.02721623553321506
$ws.Range("DH35").Value = 0.01132441717009834
$ws.Range("DN35").Value = 0.006136435233163064
$ws.Range("DO35").Value = 0.01015222382802648
$ws.Range("DP35").Value = 0.006218990444540284
$ws.Range("DQ35").Value = 0.01062805098243931
$ws.Range("DR35").Value = 0.01117235573081273
$ws.Range("DS35").Value = 0.01956583056277334
$ws.Range("DT35").Value = 0.01546572608668371
$ws.Range("DU35").Value = 0.02721623553321507
$ws.Range("DV35").Value = 0.002283504760641688
$ws.Range("DW35").Value = 0.02762348199359057
$ws.Range("DX35").Value = 8.001545544012693
$ws.Range("DY35").Value = 0.002283504760641688
$ws.Range("DZ35").Value = 11.95360915810232
$ws.Range("EA35").Value = 0.02091336478802399
$ws.Range("EB35").Value = 4.197717304840514
$ws.Range("EC35").Value = 0.0571493397041708
$ws.Range("ED35").Value = 0.04182672957604781
$ws.Range("EE35").Value = 0.01357069184865537
$ws.Range("EF35").Value = 0.001076959862540452
$ws.Range("EG35").Value = 14.9839611216823
$ws.Range("EH35").Value = 0.01357069184875303
$ws.Range("EI35").Value = 24.02099051622806
$ws.Range("EJ35").Value = 0.008225730361300703
$ws.Range("EL35").Value = 0.6404178806474745
$ws.Range("EM35").Value = 0.05056975380368477
$ws.Range("EN35").Value = 0.01645146072260143
$ws.Range("ER35").Value = 0.04580051257663533
$ws.Range("ES35").Value = 5.197001815150445
$ws.Range("ET35").Value = 0.02943510275021679
$ws.Range("FD35").Value = 0.04218008382950799
$ws.Range("FE35").Value = 0.08585713793230346
$ws.Range("FF35").Value = 0.0542254769223045
$ws.Range("FG35").Value = 0.0003154675082754932
$ws.Range("FH35").Value = 0.01944947691118454
$ws.Range("FI35").Value = 0.02253728427268023
$ws.Range("FJ35").Value = 0.01978721196251378
$ws.Range("FK35").Value = 0.00008142419552503738
$ws.Range("FP35").Value = 405.5783748941982
$ws.Range("FQ35").Value = 0.06371396799902426
$ws.Range("FR35").Value = 5.037220027622929
$ws.Range("FS35").Value = 0.3922887935245519
$ws.Range("FT35").Value = 0.005015130573029707
$ws.Range("FU35").Value = 0.002027298925603145
$ws.Range("FV35").Value = 405.5783748941982
$ws.Range("FW35").Value = 0.0000005367771389076341
$ws.Range("FX35").Value = 0.0000005627320928767214
$ws.Range("FY35").Value = 0.0000005108239979813194
$ws.Range("CZ36").Value = 0.001014312145094404
$ws.Range("DA36").Value = 0.0001151461049604975
$ws.Range("DB36").Value = 0.0001374762449219746
$ws.Range("DC36").Value = 0.000001778380627605027
$ws.Range("DD36").Value = 0.002559841961813934
$ws.Range("DE36").Value = 0.0003469518354883165
$ws.Range("DF36").Value = 0.000002298817469084363
$ws.Range("DG36").Value = 0.000348139476169094
$ws.Range("DH36").Value = 0.01033968524226356
$ws.Range("DN36").Value = 0.001014312145094404
$ws.Range("DO36").Value = 0.0001151461049604975
$ws.Range("DP36").Value = 0.001025030043825404
$ws.Range("DQ36").Value = 0.0001290085501883448
$ws.Range("DR36").Value = 0.001861516662106801
$ws.Range("DS36").Value = 0.0002475726428647713
$ws.Range("DT36").Value = 0.002548914297169838
$ws.Range("DU36").Value = 0.0003481394761690951
$ws.Range("DV36").Value = 0.005602040042932474
$ws.Range("DW36").Value = 0.0003137507959474874
$ws.Range("DX36").Value = 3.343633100803501
$ws.Range("DY36").Value = 0.005602040042932474
$ws.Range("DZ36").Value = 8.206911837435381
$ws.Range("EA36").Value = 0.01456802299473068
$ws.Range("EB36").Value = 0.09922285747058783
$ws.Range("EC36").Value = 0.08271342171097108
$ws.Range("ED36").Value = 0.02913604598946133
$ws.Range("EE36").Value = 0.007123032486102711
$ws.Range("EF36").Value = 0.0001984312494927369
$ws.Range("EG36").Value = 3.625762762679198
$ws.Range("EH36").Value = 0.007123032486102711
$ws.Range("EI36").Value = 9.401261300719026
$ws.Range("EJ36").Value = 0.009240018416360966
$ws.Range("EL36").Value = 0.08691095017693695
$ws.Range("EM36").Value = 0.01731597385631194
$ws.Range("EN36").Value = 0.01848003683272193
$ws.Range("FD36").Value = 0.06954032858510313
$ws.Range("FE36").Value = 0.05226937891196855
$ws.Range("FF36").Value = 0.08099154491702011
$ws.Range("FG36").Value = 0.0005246863693183666
$ws.Range("FH36").Value = 0.01205165470522441
$ws.Range("FI36").Value = 0.009343141726646262
$ws.Range("FJ36").Value = 0.01327880889627069
$ws.Range("FK36").Value = 0.000007772924674921477
$ws.Range("FP36").Value = 221.3879085691137
$ws.Range("FV36").Value = 221.3879085691137
$ws.Range("FW36").Value = 0.0000002930037090473629
$ws.Range("FX36").Value = 0.0000003153436431989668
$ws.Range("FY36").Value = 0.0000002706642011849274
$ws.Range("D37").Value = 0.0106019541895597
$ws.Range("T37").Value = 0.0106019541895597
$ws.Range("CZ37").Value = 0.0156668045592022
$ws.Range("DA37").Value = 0.005977635843922903
$ws.Range("DB37").Value = 0.006405107612419802
$ws.Range("DC37").Value = 0.00002007895222456721
$ws.Range("DD37").Value = 0.03953866066992996
$ws.Range("DE37").Value = 0.01616471153928484
$ws.Range("DF37").Value = 0.000002788866365078442
$ws.Range("DG37").Value = 0.0162063285868849
$ws.Range("DH37").Value = 0.01556219194656709
$ws.Range("DN37").Value = 0.0156668045592022
$ws.Range("DO37").Value = 0.005977635843922903
$ws.Range("DP37").Value = 0.01579508278003064
$ws.Range("DQ37").Value = 0.006235646658642434
$ws.Range("DR37").Value = 0.02895428398217595
$ws.Range("DS37").Value = 0.01202408384983533
$ws.Range("DT37").Value = 0.03963276587703544
$ws.Range("DU37").Value = 0.0162063285868849
$ws.Range("DV37").Value = 0.03575386247148448
$ws.Range("DW37").Value = 0.01494975892821317
$ws.Range("DX37").Value = 47.98630447190646
$ws.Range("DY37").Value = 0.03575386247148448
$ws.Range("DZ37").Value = 75.75037279329551
$ws.Range("EA37").Value = 0.05072785856291934
$ws.Range("EB37").Value = 0.3045364153961654
$ws.Range("EC37").Value = 0.250195160544648
$ws.Range("ED37").Value = 0.1014557171258382
$ws.Range("EE37").Value = 0.004968143322405316
$ws.Range("EF37").Value = 0.00639282825511926
$ws.Range("EG37").Value = 88.45710777317068
$ws.Range("EH37").Value = 0.004968143322405316
$ws.Range("EI37").Value = 115.1729835503215
$ws.Range("EJ37").Value = 0.01379522877860202
$ws.Range("EL37").Value = 0.2957444026540684
$ws.Range("EM37").Value = 0.1071207651758508
$ws.Range("EN37").Value = 0.02759045755720392
$ws.Range("FD37").Value = 0.05881854477052132
$ws.Range("FE37").Value = 0.06060349065490931
$ws.Range("FF37").Value = 0.04918240296931495
$ws.Range("FG37").Value = 0.008446571514539683
$ws.Range("FH37").Value = 0.01317604629290899
$ws.Range("FI37").Value = 0.01787119981444951
$ws.Range("FJ37").Value = 0.01411275106464393
$ws.Range("FK37").Value = 0.005655278691033835
$ws.Range("FP37").Value = 223.8619704598829
$ws.Range("FV37").Value = 223.8619704598829
$ws.Range("FW37").Value = 0.0000002962780944080793
$ws.Range("FX37").Value = 0.0000003226300966598488
$ws.Range("FY37").Value = 0.0000002699265006804926
